$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers formatted with "." as both a
# thousands separator and a decimal point (e.g. "41.222.69"), so the
# column must stay text-typed; pre-format as Text so values such as
# "3.00" or "0.783" are not reinterpreted as numbers and trimmed.
$dCells = @("D2","D3","D4","D5","D6","D9","D10","D11","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D41","D43","D44","D47","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.222.69"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "2.463.35"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "311.60"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "93.64"
$ws.Range("E6").Value = "  -6.55%  "
$ws.Range("E7").Value = "  -2.95%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -4.87%  "
$ws.Range("D10").Value = "33.23"
$ws.Range("D11").Value = "0.0777"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("E13").Value = "  -5.24%  "
$ws.Range("D14").Value = "2.842.62"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").Value = "2.444.63"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "14.89"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("D18").Value = "41.167.70"
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  -6.66%  "
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "11.23"
$ws.Range("E21").Value = "  -8.38%  "
$ws.Range("D22").Value = "68.25"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "235.37"
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").Value = "2.76"
$ws.Range("E24").Value = "  -4.00%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -6.31%  "
$ws.Range("D27").Value = "23.98"
$ws.Range("E27").Value = "  -5.88%  "
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -5.64%  "
$ws.Range("D29").Value = "9.59"
$ws.Range("E29").Value = "  -5.69%  "
$ws.Range("D30").Value = "36.22"
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("D31").Value = "152.68"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("D32").Value = "5.47"
$ws.Range("E32").Value = "  -5.27%  "
$ws.Range("D33").Value = "2.66"
$ws.Range("E33").Value = "  -5.60%  "
$ws.Range("D34").Value = "2.54"
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("D35").Value = "0.0741"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("D36").Value = "3.00"
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("D38").Value = "16.90"
$ws.Range("E38").Value = "  -7.70%  "
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("E40").Value = "  -8.36%  "
$ws.Range("D41").Value = "4.21"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D43").Value = "20.11"
$ws.Range("E43").Value = "  -10.78%  "
$ws.Range("D44").Value = "1.980.29"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  -5.43%  "
$ws.Range("E46").Value = "  -7.98%  "
$ws.Range("D47").Value = "8.67"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("D49").Value = "96.76"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("D50").Value = "74.02"
$ws.Range("E50").Value = "  -6.90%  "
$ws.Range("E51").Value = "  -7.00%  "
